$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Elements")

# Clear the "Condition(s)" values for the CodeableConcept.extension, .coding and .text rows
$ws.Range("AI4").Value = ""
$ws.Range("AI5").Value = ""
$ws.Range("AI6").Value = ""

# Fix the "Mapping: RIM Mapping" value for CodeableConcept.extension from "N/A" to "n/a"
$ws.Range("AL4").Value = "n/a"
